$wb = $excel.ActiveWorkbook

# --- "borders" sheet: rewrite the 2x2 border patterns so the first ---
# --- edge listed is the UP border (v0.16 fix), and make the sheet  ---
# --- the active one in the workbook.                               ---
$borders = $wb.Worksheets.Item("borders")

$borders.Range("B2").Value = "W _ _`nW W _"
$borders.Range("C2").Value = "W _ _`nW _ _"
$borders.Range("C3").Value = "_ W _`nW W _"
$borders.Range("D3").Value = "_ W _`n_ W W"

# New (empty) row below the existing data, taking border thickness into
# account when sizing the maze image.
$borders.Rows.Item(4).RowHeight = 52.2

# Make "borders" the active sheet / tab, with C4 selected.
$null = $borders.Activate()
$null = $borders.Range("C4").Select()
